$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.881.33'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '2.303.65'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''302.03'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '''97.28'
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("D7").Value = '''0.502'
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").Value = '''33.85'
$ws.Range("E10").Value = '  -2.39%  '
$ws.Range("D11").Value = '''0.0793'
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").Value = '''49.26'
$ws.Range("E12").Value = '  -2.60%  '
$ws.Range("E13").Value = '  +2.32%  '
$ws.Range("D14").Value = '''16.57'
$ws.Range("E14").Value = '  +7.58%  '
$ws.Range("D15").Value = '''6.78'
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").Value = '2.655.42'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '2.312.70'
$ws.Range("E17").Value = '  +1.94%  '
$ws.Range("D18").Value = '''0.798'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("D19").Value = '42.681.41'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '''11.70'
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("D21").Value = '0.0₃0900'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = '''6.03'
$ws.Range("E22").Value = '  +0.57%  '
$ws.Range("D23").Value = '''66.95'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("D24").Value = '''236.97'
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("D25").Value = '''1.99'
$ws.Range("E25").Value = '  +2.45%  '
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("D28").Value = '''24.15'
$ws.Range("E28").Value = '  -1.73%  '
$ws.Range("D29").Value = '''2.18'
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("D30").Value = '''34.62'
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("D31").Value = '''167.45'
$ws.Range("E31").Value = '  +2.05%  '
$ws.Range("D32").Value = '''9.18'
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").Value = '''4.79'
$ws.Range("E34").Value = '  +7.86%  '
$ws.Range("D35").Value = '''4.99'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = '''17.22'
$ws.Range("E36").Value = '  +3.82%  '
$ws.Range("E37").Value = '  -2.19%  '
$ws.Range("D38").Value = '''0.0699'
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("D39").Value = '''2.83'
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("D40").Value = '''0.1000'
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").Value = '''1.76'
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").Value = '''2.41'
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("D44").Value = '1.965.68'
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("D46").Value = '''17.60'
$ws.Range("E46").Value = '  -3.90%  '
$ws.Range("D47").Value = '''9.78'
$ws.Range("E47").Value = '  -3.34%  '
$ws.Range("D48").Value = '''2.85'
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").Value = '2.522.53'
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = '''53.08'
$ws.Range("E50").Value = '  -3.41%  '
$ws.Range("D51").Value = '''1.51'
$ws.Range("E51").Value = '  +1.07%  '
